$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 189, pushing the existing rows 189-204
# down to 191-206 (matches the diff: dimension grows from R204 to R206).
$ws.Rows.Item(189).Insert()
$ws.Rows.Item(189).Insert()

# New row 189
$ws.Cells.Item(189, 1).Value = 4
$ws.Cells.Item(189, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(189, 3).Value = "Los Lagos"
$ws.Cells.Item(189, 4).Value = 44516
$ws.Cells.Item(189, 5).Value = 10
$ws.Cells.Item(189, 6).Value = 100112045
$ws.Cells.Item(189, 7).Value = "Zapallo"
$ws.Cells.Item(189, 8).Value = "Camote"
$ws.Cells.Item(189, 9).Value = "Segunda"
$ws.Cells.Item(189, 10).Value = 300
$ws.Cells.Item(189, 11).Value = 600
$ws.Cells.Item(189, 12).Value = 600
$ws.Cells.Item(189, 13).Value = 600
$ws.Cells.Item(189, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(189, 15).Value = "Perú"
$ws.Cells.Item(189, 16).Value = 600
$ws.Cells.Item(189, 17).Value = 1
$ws.Cells.Item(189, 18).Value = "Hortaliza"

# New row 190
$ws.Cells.Item(190, 1).Value = 4
$ws.Cells.Item(190, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(190, 3).Value = "Los Lagos"
$ws.Cells.Item(190, 4).Value = 44516
$ws.Cells.Item(190, 5).Value = 10
$ws.Cells.Item(190, 6).Value = 100112045
$ws.Cells.Item(190, 7).Value = "Zapallo"
$ws.Cells.Item(190, 8).Value = "Paine"
$ws.Cells.Item(190, 9).Value = "1a (guarda)"
$ws.Cells.Item(190, 10).Value = 1100
$ws.Cells.Item(190, 11).Value = 400
$ws.Cells.Item(190, 12).Value = 400
$ws.Cells.Item(190, 13).Value = 400
$ws.Cells.Item(190, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(190, 15).Value = "Región Metropolitana"
$ws.Cells.Item(190, 16).Value = 400
$ws.Cells.Item(190, 17).Value = 1
$ws.Cells.Item(190, 18).Value = "Hortaliza"
